# Planning_TPI_Forestier.xlsx - "Première version qui marche bien"
#
# This script reproduces the edits captured in the commit:
#  1. Sheet "Planning effectif" becomes the active tab (was "Tâches réalisés").
#  2. A handful of cells on "Planning effectif" get re-styled (half-day marks).
#  3. View/selection state changes on all three sheets.
#  4. New rows (14.mai) are appended to "Tâches réalisés", with three new
#     shared strings plus the "Total 14.mai" summary row/formula.

$wb = $excel.ActiveWorkbook

$wsPrev = $wb.Worksheets.Item(1)   # "Planning prévisionel"
$wsEff  = $wb.Worksheets.Item(2)   # "Planning effectif"
$wsReal = $wb.Worksheets.Item(3)   # "Tâches réalisés"

function CopyFormat($ws, $srcAddr, $dstAddr) {
    $src = $ws.Range($srcAddr)
    $src.Copy()
    $dst = $ws.Range($dstAddr)
    $dst.PasteSpecial(-4122)  # xlPasteFormats
}

# ---------------------------------------------------------------------------
# 1) "Planning effectif" - restyle a few half-day cells and drop L9
# ---------------------------------------------------------------------------

CopyFormat $wsEff "J7"  "K7"    # K7: 28 -> 24
CopyFormat $wsEff "J7"  "L8"    # L8: 28 -> 24
$wsEff.Range("L9").Clear()      # L9 cell removed entirely
CopyFormat $wsEff "J7"  "L11"   # L11: 28 -> 24
CopyFormat $wsEff "J7"  "L12"   # L12: 28 -> 24
CopyFormat $wsEff "J7"  "K15"   # K15: 31 -> 24
CopyFormat $wsEff "J7"  "L15"   # L15: 31 -> 24

# ---------------------------------------------------------------------------
# 2) New rows for 14.mai on "Tâches réalisés"
# ---------------------------------------------------------------------------

# Row 31 - like row 29 (Date / Tâche / Temps / Commentaire, normal height)
CopyFormat $wsReal "A29" "A31"
CopyFormat $wsReal "B29" "B31"
CopyFormat $wsReal "C29" "C31"
CopyFormat $wsReal "D29" "D31"
$wsReal.Cells.Item(31, 1).Value = 43599
$wsReal.Cells.Item(31, 2).Value = "Analyse d'image plus précise"
$wsReal.Cells.Item(31, 3).Value = 0.16666666666666666
$wsReal.Cells.Item(31, 4).Value = "L'analyse est à mon avis suffisante pour le projet. A voir avec M.Bonvin"

# Row 32 - like row 6 (Date / Tâche / Temps / Commentaire, tall wrapped row)
CopyFormat $wsReal "A6"  "A32"
CopyFormat $wsReal "B6"  "B32"
CopyFormat $wsReal "C10" "C32"
CopyFormat $wsReal "D6"  "D32"
$wsReal.Cells.Item(32, 1).Value = 43599
$wsReal.Cells.Item(32, 2).Value = "Implémentation de la télécommande du robot"
$wsReal.Cells.Item(32, 3).Value = 0.125
$wsReal.Cells.Item(32, 4).Value = "Direction, vitesse/compensation, mode manuel/auto"
$wsReal.Rows.Item(32).RowHeight = 28.8

# Row 33 - like row 26 (Date / Tâche / Temps only)
CopyFormat $wsReal "A26" "A33"
CopyFormat $wsReal "B26" "B33"
$wsReal.Cells.Item(33, 1).Value = 43599
$wsReal.Cells.Item(33, 2).Value = "Documentation"
$wsReal.Cells.Item(33, 3).Value = 0.041666666666666664

# Row 34 - like row 30 (Tâche / Temps totals row with SUM formula)
CopyFormat $wsReal "B30" "B34"
$wsReal.Cells.Item(34, 2).Value = "Total 14.mai"
$wsReal.Range("C34").Formula = "=SUM(C31:C33)"

# Apply the formats of C33/C34 *after* values/formula are in place: pasting
# formats onto a SUM cell before its formula is (re)written confuses this
# engine's calc cache and leaves a stale cached result.
CopyFormat $wsReal "C26" "C33"
CopyFormat $wsReal "C30" "C34"

$excel.Calculate()

# ---------------------------------------------------------------------------
# 3) View / selection state
# ---------------------------------------------------------------------------

$wsPrev.Activate()
$wsPrev.Range("Q10").Select()

$wsReal.Activate()
$wsReal.Range("D33").Select()

$wsEff.Activate()
$wsEff.Range("AA12").Select()
